$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update phone numbers (column J) for rows 2-34, clearing the border style ---
$ws.Range("J2").Value = 5616677351
$ws.Range("J2").Style = "Normal"
$ws.Range("J3").Value = 5614858433
$ws.Range("J3").Style = "Normal"
$ws.Range("J4").Value = 5615869888
$ws.Range("J4").Style = "Normal"
$ws.Range("J5").Value = 5614714457
$ws.Range("J5").Style = "Normal"
$ws.Range("J6").Value = 5619647774
$ws.Range("J6").Style = "Normal"
$ws.Range("J7").Value = 5614786795
$ws.Range("J7").Style = "Normal"
$ws.Range("J8").Value = 5619744199
$ws.Range("J8").Style = "Normal"
$ws.Range("J9").Value = 5618645137
$ws.Range("J9").Style = "Normal"
$ws.Range("J10").Value = 5616914758
$ws.Range("J10").Style = "Normal"
$ws.Range("J11").Value = 5616667133
$ws.Range("J11").Style = "Normal"
$ws.Range("J12").Value = 5615945674
$ws.Range("J12").Style = "Normal"
$ws.Range("J13").Value = 5615853684
$ws.Range("J13").Style = "Normal"
$ws.Range("J14").Value = 5615938591
$ws.Range("J14").Style = "Normal"
$ws.Range("J15").Value = 5616677373
$ws.Range("J15").Style = "Normal"
$ws.Range("J16").Value = 5619831869
$ws.Range("J16").Style = "Normal"
$ws.Range("J17").Value = 5618956315
$ws.Range("J17").Style = "Normal"
$ws.Range("J18").Value = 5616914654
$ws.Range("J18").Style = "Normal"
$ws.Range("J19").Value = 5618954116
$ws.Range("J19").Style = "Normal"
$ws.Range("J20").Value = 5618956315
$ws.Range("J20").Style = "Normal"
$ws.Range("J21").Value = 5613917133
$ws.Range("J21").Style = "Normal"
$ws.Range("J22").Value = 5616695198
$ws.Range("J22").Style = "Normal"
$ws.Range("J23").Value = 5615955997
$ws.Range("J23").Style = "Normal"
$ws.Range("J24").Value = 5614854416
$ws.Range("J24").Style = "Normal"
$ws.Range("J25").Value = 5613643175
$ws.Range("J25").Style = "Normal"
$ws.Range("J26").Value = 5614788353
$ws.Range("J26").Style = "Normal"
$ws.Range("J27").Value = 5616921415
$ws.Range("J27").Style = "Normal"
$ws.Range("J28").Value = 5615749433
$ws.Range("J28").Style = "Normal"
$ws.Range("J29").Value = 5619843437
$ws.Range("J29").Style = "Normal"
$ws.Range("J30").Value = 5619813471
$ws.Range("J30").Style = "Normal"
$ws.Range("J31").Value = 5614975479
$ws.Range("J31").Style = "Normal"
$ws.Range("J32").Value = 5617755983
$ws.Range("J32").Style = "Normal"
$ws.Range("J33").Value = 5616744552
$ws.Range("J33").Style = "Normal"
$ws.Range("J34").Value = 5617671652
$ws.Range("J34").Style = "Normal"

# --- Append new school records (Alborz / Karaj district 4) in rows 97-117 ---
$ws.Range("A97").Value = 'البرز'
$ws.Range("B97").Value = 'کرج ناحيه 4'
$ws.Range("C97").Value = 'شهيدان  هداوند'
$ws.Range("D97").Value = 'دوره ابتدايي توصيفي'
$ws.Range("E97").Value = 'دخترانه'
$ws.Range("F97").Value = 'دولتي'
$ws.Range("G97").Value = 'عادي'
$ws.Range("H97").Value = 'عادي'
$ws.Range("I97").Value = 4552018
$ws.Range("J97").Value = 319777333
$ws.Range("K97").Value = 'حصار ک بالا خيابان آقارضايي آموزشگاه شهيدان هداوند2'

$ws.Range("A98").Value = 'البرز'
$ws.Range("B98").Value = 'کرج ناحيه 4'
$ws.Range("C98").Value = 'شهداي کمالشهر'
$ws.Range("D98").Value = 'دوره متوسطه اول'
$ws.Range("E98").Value = 'پسرانه'
$ws.Range("F98").Value = 'دولتي'
$ws.Range("G98").Value = 'عادي'
$ws.Range("H98").Value = 'عادي'
$ws.Range("I98").Value = 4703847
$ws.Range("J98").Value = 319976577
$ws.Range("K98").Value = 'کمالشهر_خ وليعصرجنوبي _خ شهدا'

$ws.Range("A99").Value = 'البرز'
$ws.Range("B99").Value = 'کرج ناحيه 4'
$ws.Range("C99").Value = 'شهيد صياد شيرازي'
$ws.Range("D99").Value = 'دوره ابتدايي توصيفي'
$ws.Range("E99").Value = 'پسرانه'
$ws.Range("F99").Value = 'دولتي'
$ws.Range("G99").Value = 'عادي'
$ws.Range("H99").Value = 'عادي'
$ws.Range("I99").Value = 34703847
$ws.Range("J99").Value = '-'
$ws.Range("K99").Value = 'کمالشهر جنب مسجد علي ابن ابي طالب- مدرسه ش صياد شيرازي'

$ws.Range("A100").Value = 'البرز'
$ws.Range("B100").Value = 'کرج ناحيه 4'
$ws.Range("C100").Value = 'شهيدستارلطفي'
$ws.Range("D100").Value = 'متوسطه دوم - نظري'
$ws.Range("E100").Value = 'پسرانه'
$ws.Range("F100").Value = 'دولتي'
$ws.Range("G100").Value = 'عادي'
$ws.Range("H100").Value = 'عادي'
$ws.Range("I100").Value = '-'
$ws.Range("J100").Value = '-'
$ws.Range("K100").Value = 'کيانمهر'

$ws.Range("A101").Value = 'البرز'
$ws.Range("B101").Value = 'کرج ناحيه 4'
$ws.Range("C101").Value = 'هنرستان امام خميني (ره )(1)'
$ws.Range("D101").Value = 'متوسطه دوم - هنرستان فني'
$ws.Range("E101").Value = 'پسرانه'
$ws.Range("F101").Value = 'دولتي'
$ws.Range("G101").Value = 'عادي'
$ws.Range("H101").Value = 'هيات امنايي'
$ws.Range("I101").Value = 3528485
$ws.Range("J101").Value = 318383464
$ws.Range("K101").Value = 'فاز4مهرشهر_بلوار گلها _چهارراه هنرستان خ 406شرقي-پ192'

$ws.Range("A102").Value = 'البرز'
$ws.Range("B102").Value = 'کرج ناحيه 4'
$ws.Range("C102").Value = 'حضرت مريم (1)'
$ws.Range("D102").Value = 'متوسطه دوم - نظري'
$ws.Range("E102").Value = 'دخترانه'
$ws.Range("F102").Value = 'دولتي'
$ws.Range("G102").Value = 'عادي'
$ws.Range("H102").Value = 'عادي'
$ws.Range("I102").Value = 3408607
$ws.Range("J102").Value = 318579855
$ws.Range("K102").Value = 'مهرشهر_بلوارارم _بلواردانش _خ 100_خ مريم'

$ws.Range("A103").Value = 'البرز'
$ws.Range("B103").Value = 'کرج ناحيه 4'
$ws.Range("C103").Value = 'شهيد شهسواري(2)'
$ws.Range("D103").Value = 'دوره متوسطه اول'
$ws.Range("E103").Value = 'پسرانه'
$ws.Range("F103").Value = 'دولتي'
$ws.Range("G103").Value = 'عادي'
$ws.Range("H103").Value = 'عادي'
$ws.Range("I103").Value = 34662060
$ws.Range("J103").Value = '-'
$ws.Range("K103").Value = 'کرج- حصارک پايين -رضاشهر-انتهاي خيابان فروردين'

$ws.Range("A104").Value = 'البرز'
$ws.Range("B104").Value = 'کرج ناحيه 4'
$ws.Range("C104").Value = 'لقمان حکيم(1)'
$ws.Range("D104").Value = 'دوره متوسطه اول'
$ws.Range("E104").Value = 'پسرانه'
$ws.Range("F104").Value = 'دولتي'
$ws.Range("G104").Value = 'عادي'
$ws.Range("H104").Value = 'عادي'
$ws.Range("I104").Value = 34801188
$ws.Range("J104").Value = '-'
$ws.Range("K104").Value = 'کرج-پيشاهنگي-گلدشت'

$ws.Range("A105").Value = 'البرز'
$ws.Range("B105").Value = 'کرج ناحيه 4'
$ws.Range("C105").Value = 'هدايت'
$ws.Range("D105").Value = 'دوره ابتدايي توصيفي'
$ws.Range("E105").Value = 'دخترانه'
$ws.Range("F105").Value = 'دولتي'
$ws.Range("G105").Value = 'عادي'
$ws.Range("H105").Value = 'عادي'
$ws.Range("I105").Value = '-'
$ws.Range("J105").Value = '-'
$ws.Range("K105").Value = '-'

$ws.Range("A106").Value = 'البرز'
$ws.Range("B106").Value = 'کرج ناحيه 4'
$ws.Range("C106").Value = 'شهيد پرورش'
$ws.Range("D106").Value = 'دوره ابتدايي توصيفي'
$ws.Range("E106").Value = 'پسرانه'
$ws.Range("F106").Value = 'دولتي'
$ws.Range("G106").Value = 'عادي'
$ws.Range("H106").Value = 'عادي'
$ws.Range("I106").Value = 33214848
$ws.Range("J106").Value = 123454
$ws.Range("K106").Value = 'شهرک کيان مهر-خ نبرد اهواز - مدرسه شهيدمحمدپرورشي'

$ws.Range("A107").Value = 'البرز'
$ws.Range("B107").Value = 'کرج ناحيه 4'
$ws.Range("C107").Value = 'امام رضا(ع )2'
$ws.Range("D107").Value = 'دوره ابتدايي توصيفي'
$ws.Range("E107").Value = 'دخترانه'
$ws.Range("F107").Value = 'دولتي'
$ws.Range("G107").Value = 'عادي'
$ws.Range("H107").Value = 'عادي'
$ws.Range("I107").Value = 3315356
$ws.Range("J107").Value = 318695759
$ws.Range("K107").Value = 'جاده قزلحصارروبروي بي سيم شهرک سهرابيه'

$ws.Range("A108").Value = 'البرز'
$ws.Range("B108").Value = 'کرج ناحيه 4'
$ws.Range("C108").Value = 'شهيد باهنر'
$ws.Range("D108").Value = 'دوره ابتدايي توصيفي'
$ws.Range("E108").Value = 'دخترانه'
$ws.Range("F108").Value = 'دولتي'
$ws.Range("G108").Value = 'عادي'
$ws.Range("H108").Value = 'عادي'
$ws.Range("I108").Value = 33203180
$ws.Range("J108").Value = 318761748
$ws.Range("K108").Value = 'کيان مهر ميدان مهرگان بوستان هفتم'

$ws.Range("A109").Value = 'البرز'
$ws.Range("B109").Value = 'کرج ناحيه 4'
$ws.Range("C109").Value = 'وحدت اسلامي (1)'
$ws.Range("D109").Value = 'متوسطه دوم - هنرستان کاردانش'
$ws.Range("E109").Value = 'پسرانه'
$ws.Range("F109").Value = 'دولتي'
$ws.Range("G109").Value = 'عادي'
$ws.Range("H109").Value = 'عادي'
$ws.Range("I109").Value = 4553805
$ws.Range("J109").Value = 319767745
$ws.Range("K109").Value = 'حصارک بالا_روبروي مجتمع ورزشي ايثار'

$ws.Range("A110").Value = 'البرز'
$ws.Range("B110").Value = 'کرج ناحيه 4'
$ws.Range("C110").Value = 'وليعصر(عج)'
$ws.Range("D110").Value = 'دوره متوسطه اول'
$ws.Range("E110").Value = 'پسرانه'
$ws.Range("F110").Value = 'دولتي'
$ws.Range("G110").Value = 'عادي'
$ws.Range("H110").Value = 'عادي'
$ws.Range("I110").Value = 3212728
$ws.Range("J110").Value = 318761748
$ws.Range("K110").Value = 'کيانمهر_بلواراميرکبير_جنب ميدان امام خميني (ره )'

$ws.Range("A111").Value = 'البرز'
$ws.Range("B111").Value = 'کرج ناحيه 4'
$ws.Range("C111").Value = 'صداقت'
$ws.Range("D111").Value = 'دوره ابتدايي توصيفي'
$ws.Range("E111").Value = 'پسرانه'
$ws.Range("F111").Value = 'غيردولتي'
$ws.Range("G111").Value = 'عادي'
$ws.Range("H111").Value = 'عادي'
$ws.Range("I111").Value = 34801130
$ws.Range("J111").Value = '-'
$ws.Range("K111").Value = 'کرج-خرمدشت-ميثم يک-بهار اول-پلاک37'

$ws.Range("A112").Value = 'البرز'
$ws.Range("B112").Value = 'کرج ناحيه 4'
$ws.Range("C112").Value = 'فرازين'
$ws.Range("D112").Value = 'دوره ابتدايي توصيفي'
$ws.Range("E112").Value = 'دخترانه'
$ws.Range("F112").Value = 'غيردولتي'
$ws.Range("G112").Value = 'عادي'
$ws.Range("H112").Value = 'عادي'
$ws.Range("I112").Value = 33420115
$ws.Range("J112").Value = '-'
$ws.Range("K112").Value = 'کرج-فاز2مهرشهر-بلوار شهرداري-خيابان202-پلاک382/1 پيش و ابتدايي فرازين'

$ws.Range("A113").Value = 'البرز'
$ws.Range("B113").Value = 'کرج ناحيه 4'
$ws.Range("C113").Value = 'حضرت امير(ع)'
$ws.Range("D113").Value = 'دوره متوسطه اول'
$ws.Range("E113").Value = 'پسرانه'
$ws.Range("F113").Value = 'غيردولتي'
$ws.Range("G113").Value = 'عادي'
$ws.Range("H113").Value = 'عادي'
$ws.Range("I113").Value = 3509123
$ws.Range("J113").Value = 313965359
$ws.Range("K113").Value = 'کرج-خيابان درختي-نرسيده به سه راه تهران-پلاک325-متوسطه دوره اول حضرت امير(ع)'

$ws.Range("A114").Value = 'البرز'
$ws.Range("B114").Value = 'کرج ناحيه 4'
$ws.Range("C114").Value = 'شادان'
$ws.Range("D114").Value = 'دوره متوسطه اول'
$ws.Range("E114").Value = 'دخترانه'
$ws.Range("F114").Value = 'غيردولتي'
$ws.Range("G114").Value = 'عادي'
$ws.Range("H114").Value = 'عادي'
$ws.Range("I114").Value = 34613059
$ws.Range("J114").Value = '-'
$ws.Range("K114").Value = 'کرج-حصارک - خيابان برزنت-90دستگاه اول-پلاک34-متوسطه اول شادان'

$ws.Range("A115").Value = 'البرز'
$ws.Range("B115").Value = 'کرج ناحيه 4'
$ws.Range("C115").Value = 'عصر تلاش'
$ws.Range("D115").Value = 'دوره ابتدايي توصيفي'
$ws.Range("E115").Value = 'پسرانه'
$ws.Range("F115").Value = 'غيردولتي'
$ws.Range("G115").Value = 'عادي'
$ws.Range("H115").Value = 'عادي'
$ws.Range("I115").Value = 34516780
$ws.Range("J115").Value = '-'
$ws.Range("K115").Value = '-'

$ws.Range("A116").Value = 'البرز'
$ws.Range("B116").Value = 'کرج ناحيه 4'
$ws.Range("C116").Value = 'نيوشا'
$ws.Range("D116").Value = 'دوره متوسطه اول'
$ws.Range("E116").Value = 'دخترانه'
$ws.Range("F116").Value = 'غيردولتي'
$ws.Range("G116").Value = 'عادي'
$ws.Range("H116").Value = 'عادي'
$ws.Range("I116").Value = 33513094
$ws.Range("J116").Value = '-'
$ws.Range("K116").Value = 'کرج-خيايان45متري گلشهر-کوچه مينا-پلاک35-متوسطه اول نيوشا'

$ws.Range("A117").Value = 'البرز'
$ws.Range("B117").Value = 'کرج ناحيه 4'
$ws.Range("C117").Value = 'انديشه'
$ws.Range("D117").Value = 'متوسطه دوم - هنرستان کاردانش'
$ws.Range("E117").Value = 'پسرانه'
$ws.Range("F117").Value = 'غيردولتي'
$ws.Range("G117").Value = 'عادي'
$ws.Range("H117").Value = 'عادي'
$ws.Range("I117").Value = 4641296
$ws.Range("J117").Value = 313891491
$ws.Range("K117").Value = 'کرج-45متري گلشهر-آذرشرقي-پلاک14-کاردانش انديشه'

# --- Copy the formatting (borders/fill/alignment) of the last pre-existing row down onto the new rows ---
$ws.Range("A96:K96").Copy()
$ws.Range("A97:K117").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Restore view state: scroll position + active selection ---
$win = $excel.ActiveWindow
$win.ScrollRow = 25
$win.ScrollColumn = 1
$ws.Range("J35").Select()

